$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (DAMSLTag, DialogAct) for columns I and J
$updates = @(
    @{Row=3;   I="sv"; J="Statement-opinion"},
    @{Row=32;  I="sv"; J="Statement-opinion"},
    @{Row=39;  I="aa"; J="Agree/Accept"},
    @{Row=45;  I="sd"; J="Statement-non-opinion"},
    @{Row=60;  I="qy"; J="Yes-No-Question"},
    @{Row=61;  I="aa"; J="Agree/Accept"},
    @{Row=65;  I="b";  J="Acknowledge (Backchannel)"},
    @{Row=66;  I="sd"; J="Statement-non-opinion"},
    @{Row=84;  I="sd"; J="Statement-non-opinion"},
    @{Row=101; I="ba"; J="Appreciation"}
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.I
    $ws.Range("J" + $u.Row).Value = $u.J
}
